# Update cryptocurrency price (D) and 1h volume change (E) columns
# per the Nov 4 2023 02:25:57 UTC GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.104.35"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "1.856.60"
$ws.Range("E3").Value = "  +2.91%  "

$ws.Range("E4").Value = "  +0.57%  "

$ws.Range("D5").Value = "'233.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.85%  "

$ws.Range("D6").Value = "'0.620"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.99%  "

$ws.Range("E7").Value = "  +0.50%  "

$ws.Range("D8").Value = "'40.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.41%  "

$ws.Range("D9").Value = "'0.333"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.40%  "

$ws.Range("D10").Value = "'0.0695"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.61%  "

$ws.Range("D11").Value = "'0.0986"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "2.121.25"
$ws.Range("E12").Value = "  +2.78%  "

$ws.Range("D13").Value = "'11.49"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.44%  "

$ws.Range("D14").Value = "1.855.53"
$ws.Range("E14").Value = "  +2.95%  "

$ws.Range("D15").Value = "'0.677"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.51%  "

$ws.Range("D16").Value = "'4.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.28%  "

$ws.Range("D17").Value = "35.146.52"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("D18").Value = "'70.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.06%  "

$ws.Range("D19").Value = "0.0₃0795"
$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("D20").Value = "'241.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "

$ws.Range("D21").Value = "'12.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.50%  "

$ws.Range("D22").Value = "'4.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.83%  "

$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").Value = "'2.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.30%  "

$ws.Range("D25").Value = "'173.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").Value = "'7.88"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'17.62"
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = "  +4.73%  "

$ws.Range("D29").Value = "'1.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.36%  "

$ws.Range("E30").Value = "  +0.43%  "

$ws.Range("D31").Value = "'0.0557"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.92%  "

$ws.Range("D32").Value = "'3.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").Value = "'4.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.69%  "

$ws.Range("D34").Value = "'1.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +24.01%  "

$ws.Range("D35").Value = "'1.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.85%  "

$ws.Range("D36").Value = "'0.765"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.60%  "

$ws.Range("E37").Value = "  +8.14%  "

$ws.Range("D38").Value = "'1.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.89%  "

$ws.Range("D39").Value = "'90.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("D40").Value = "1.356.53"
$ws.Range("E40").Value = "  +4.12%  "

$ws.Range("E41").Value = "  +3.56%  "

$ws.Range("D42").Value = "'14.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.95%  "

$ws.Range("D43").Value = "'2.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.52%  "

$ws.Range("E44").Value = "  -1.41%  "

$ws.Range("D45").Value = "'2.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.89%  "

$ws.Range("D46").Value = "'0.0531"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.71%  "

$ws.Range("D47").Value = "'6.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.54%  "

$ws.Range("D48").Value = "2.037.99"
$ws.Range("E48").Value = "  +2.75%  "

$ws.Range("D49").Value = "'3.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +20.88%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("E51").Value = "  +0.43%  "

